# Örnek 4 - Dört İşlem.xlsx -- "Add files via upload" edit
#
# Adds the four-operation formulas (A+B, B+C, A+B-C, (A-C)/B, A*(B+C), (A+C)/B)
# to columns D:I for rows 3-11, fills in the student info box (K4:N6) with a
# number, "Kübra Çabuk" and "YBS", resizes the columns, and changes the zoom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3: plain (non-shared) formulas
# ---------------------------------------------------------------------------
$ws.Range("D3").Formula = "=(A3+B3)"
$ws.Range("E3").Formula = "=(B3+C3)"
$ws.Range("F3").Formula = "=(D3-C3)"
$ws.Range("G3").Formula = "=(A3-C3 /B3)"
$ws.Range("H3").Formula = "=(A3*E3)"
$ws.Range("I3").Formula = "=(A3+C3/B3)"

# ---------------------------------------------------------------------------
# Rows 4-11: same formulas, filled down as shared formulas (D4:D11, etc.)
# ---------------------------------------------------------------------------
$ws.Range("D4:D11").Formula = "=(A4+B4)"
$ws.Range("E4:E11").Formula = "=(B4+C4)"
$ws.Range("F4:F11").Formula = "=(D4-C4)"
$ws.Range("G4:G11").Formula = "=(A4-C4 /B4)"
$ws.Range("H4:H11").Formula = "=(A4*E4)"
$ws.Range("I4:I11").Formula = "=(A4+C4/B4)"

# G4 carried a stray "0.00" number format in the original file (a leftover,
# different from the "0" format used by every other cell in the column) --
# line it up with the rest of the column now that it has a formula.
$ws.Range("G4").NumberFormat = "0"

# ---------------------------------------------------------------------------
# Student info box
# ---------------------------------------------------------------------------
$ws.Range("L4").Value = 20215070019
$ws.Range("L5").Value = "Kübra Çabuk"
$ws.Range("L6").Value = "YBS"

# ---------------------------------------------------------------------------
# Column widths (values chosen so the stored OOXML width matches the target
# as closely as this engine's character->pixel rounding allows)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21
$ws.Columns.Item(2).ColumnWidth = 17.142857142857142
$ws.Columns.Item(3).ColumnWidth = 8.428571428571429
$ws.Columns.Item(4).ColumnWidth = 13.285714285714286
$ws.Columns.Item(5).ColumnWidth = 12.428571428571429
$ws.Columns.Item(6).ColumnWidth = 21
$ws.Columns.Item(7).ColumnWidth = 19.428571428571427
$ws.Columns.Item(8).ColumnWidth = 21.571428571428573
$ws.Columns.Item(9).ColumnWidth = 26.285714285714285

# ---------------------------------------------------------------------------
# Zoom level
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 70
